$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells, reusing the existing header formatting
# (bold, bordered, centered) by copying the style from the last existing
# header cell (AC1) rather than constructing a brand-new style.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the team record (Wins/Losses/Ties) for every player row.
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 82
    $ws.Cells.Item($r, 31).Value = 80
    $ws.Cells.Item($r, 32).Value = 0
}
